$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Nous en avons identifié quatre dont un état composite,
# qui contient deux autres états." -> "... qui contient trois autres
# états." with "trois" becoming its own run and a "_GoBack" bookmark
# landing right after it (mirrors Word's own "select & retype" +
# last-edit-position bookmark behaviour).
# -----------------------------------------------------------------

# 1a. Plain text substitution first (deux -> trois).
$r1 = $d.Content
$r1.Find.Execute("deux autres états.", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "trois autres états.", 2) | Out-Null

# 1b. Locate the freshly inserted word "trois" so we can split the run
# both just before and just after it.
$triois = $d.Content
$triois.Find.Execute("trois autres") | Out-Null
$troisStart = $triois.Start
$troisEnd = $troisStart + 5   # length of "trois"

# Force a run boundary right before "trois" using a throw-away bookmark
# (adding a bookmark splits the enclosing run; deleting the bookmark
# afterwards leaves the split in place).
$splitBefore = $d.Range($troisStart, $troisStart)
$d.Bookmarks.Add("zzTempSplit", $splitBefore) | Out-Null
$d.Bookmarks.Item("zzTempSplit").Delete()

# Put the real "_GoBack" bookmark right after "trois" - this both
# splits the run there AND records the bookmark, matching the diff.
$gobackRange = $d.Range($troisEnd, $troisEnd)
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null

# -----------------------------------------------------------------
# Change 2: the old "_GoBack" bookmark used to sit inside "transitions"
# (splitting "Nous pouvons voir que les t" / "ransitions, ... sur ").
# Adding the new "_GoBack" bookmark above already removed it (a
# document can only have one "_GoBack"), but the two runs it used to
# separate are still split. Re-merge them into a single run, exactly
# like the target XML.
# -----------------------------------------------------------------

$merge = $d.Content
$merge.Find.Execute("Nous pouvons voir que les transitions, pour la plupart, portent des événements de type change, sur ") | Out-Null
$mergeText = $merge.Text
$mergeRange = $d.Range($merge.Start, $merge.End)
$mergeRange.Delete()
$reinsert = $d.Range($merge.Start, $merge.Start)
$reinsert.InsertBefore($mergeText)

Write-Output "done"
